# "Generate Report for handback"
#
# The localization-status report is regenerated after a handback: the
# "Ready for handoff" status becomes "Handed back: in sync with en-US" for
# every tracked source file, each zh-cn/de-de detail sheet grows two new
# columns (E = Latest Target File, F = Latest Handback File) that mirror the
# source .md / handoff .xlf hyperlinks, and the "Latest Handback DateTime"
# column (G) is stamped with the handback time instead of the empty-default
# "0001-01-01 00:00:00".

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

# --- Overview sheet: just the status text, in both language columns -------
$overview.Range("B2").Value = $newStatus
$overview.Range("C2").Value = $newStatus
$overview.Range("B3").Value = $newStatus
$overview.Range("C3").Value = $newStatus

# --- Shared hyperlink target addresses (reused for the new columns) -------
$mdUrl1   = "https://github.com/OpenLocalizationTest/oltest/blob/38ee9270e89bef23a4a969f8c9a191c76a56a004/e2e/2f896ebd-0648-4060-aacb-62692ef7c544.md"
$mdUrl2   = "https://github.com/OpenLocalizationTest/oltest/blob/38ee9270e89bef23a4a969f8c9a191c76a56a004/e2e/6329960e-3640-4a37-8321-56d9dcf5ee52.md"

$zhXlfUrl1 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7a027417ba2e653496eddd2c5b7512451b3fb4c6/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/2f896ebd-0648-4060-aacb-62692ef7c544.9e4bc652c0f3b7031a37aa768eb56a3aab3d784e.zh-cn.xlf"
$zhXlfUrl2 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7a027417ba2e653496eddd2c5b7512451b3fb4c6/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/6329960e-3640-4a37-8321-56d9dcf5ee52.2daef6e2b6c383f8d79ffb9a56a7f70dbdb68b3b.zh-cn.xlf"

$deXlfUrl1 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c890edf12701571864038d9d3f09c8be710000be/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/2f896ebd-0648-4060-aacb-62692ef7c544.9e4bc652c0f3b7031a37aa768eb56a3aab3d784e.de-de.xlf"
$deXlfUrl2 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c890edf12701571864038d9d3f09c8be710000be/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/6329960e-3640-4a37-8321-56d9dcf5ee52.2daef6e2b6c383f8d79ffb9a56a7f70dbdb68b3b.de-de.xlf"

$mdName1  = "2f896ebd-0648-4060-aacb-62692ef7c544.md"
$mdName2  = "6329960e-3640-4a37-8321-56d9dcf5ee52.md"

# --- zh-cn detail sheet -----------------------------------------------------
$zhcn.Range("B2").Value = $newStatus
$zhcn.Range("B3").Value = $newStatus

$zhcn.Hyperlinks.Add($zhcn.Range("E2"), $mdUrl1, $null, $null, $mdName1)
$zhcn.Hyperlinks.Add($zhcn.Range("F2"), $zhXlfUrl1, $null, $null, "2f896ebd-0648-4060-aacb-62692ef7c544.9e4bc652c0f3b7031a37aa768eb56a3aab3d784e.zh-cn.xlf")
$zhcn.Hyperlinks.Add($zhcn.Range("E3"), $mdUrl2, $null, $null, $mdName2)
$zhcn.Hyperlinks.Add($zhcn.Range("F3"), $zhXlfUrl2, $null, $null, "6329960e-3640-4a37-8321-56d9dcf5ee52.2daef6e2b6c383f8d79ffb9a56a7f70dbdb68b3b.zh-cn.xlf")

$zhcn.Range("G2").Value = "2016-01-25 09:19:23"
$zhcn.Range("G3").Value = "2016-01-25 09:19:23"

# --- de-de detail sheet -----------------------------------------------------
$dede.Range("B2").Value = $newStatus
$dede.Range("B3").Value = $newStatus

$dede.Hyperlinks.Add($dede.Range("E2"), $mdUrl1, $null, $null, $mdName1)
$dede.Hyperlinks.Add($dede.Range("F2"), $deXlfUrl1, $null, $null, "2f896ebd-0648-4060-aacb-62692ef7c544.9e4bc652c0f3b7031a37aa768eb56a3aab3d784e.de-de.xlf")
$dede.Hyperlinks.Add($dede.Range("E3"), $mdUrl2, $null, $null, $mdName2)
$dede.Hyperlinks.Add($dede.Range("F3"), $deXlfUrl2, $null, $null, "6329960e-3640-4a37-8321-56d9dcf5ee52.2daef6e2b6c383f8d79ffb9a56a7f70dbdb68b3b.de-de.xlf")

$dede.Range("G2").Value = "2016-01-25 09:19:45"
$dede.Range("G3").Value = "2016-01-25 09:19:45"

Write-Host "Report regenerated for handback"
